# Update scraped_at timestamps on the "snapshot" sheet (column K, rows 2-48)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("snapshot")

$ws1.Range("K2").Value = "2025-11-03T10:55:08.594127+00:00"
$ws1.Range("K3").Value = "2025-11-03T10:55:10.750751+00:00"
$ws1.Range("K4").Value = "2025-11-03T10:55:10.750773+00:00"
$ws1.Range("K5").Value = "2025-11-03T10:55:10.750781+00:00"
$ws1.Range("K6").Value = "2025-11-03T10:55:10.750788+00:00"
$ws1.Range("K7").Value = "2025-11-03T10:55:13.334977+00:00"
$ws1.Range("K8").Value = "2025-11-03T10:55:13.335009+00:00"
$ws1.Range("K9").Value = "2025-11-03T10:55:13.335030+00:00"
$ws1.Range("K10").Value = "2025-11-03T10:55:15.906289+00:00"
$ws1.Range("K11").Value = "2025-11-03T10:55:18.068454+00:00"
$ws1.Range("K12").Value = "2025-11-03T10:55:18.068486+00:00"
$ws1.Range("K13").Value = "2025-11-03T10:55:18.068509+00:00"
$ws1.Range("K14").Value = "2025-11-03T10:55:20.424177+00:00"
$ws1.Range("K15").Value = "2025-11-03T10:55:20.424207+00:00"
$ws1.Range("K16").Value = "2025-11-03T10:55:20.424227+00:00"
$ws1.Range("K17").Value = "2025-11-03T10:55:27.741224+00:00"
$ws1.Range("K18").Value = "2025-11-03T10:55:29.834442+00:00"
$ws1.Range("K19").Value = "2025-11-03T10:55:31.955401+00:00"
$ws1.Range("K20").Value = "2025-11-03T10:55:34.117199+00:00"
$ws1.Range("K21").Value = "2025-11-03T10:55:34.117231+00:00"
$ws1.Range("K22").Value = "2025-11-03T10:55:34.117249+00:00"
$ws1.Range("K23").Value = "2025-11-03T10:55:36.701819+00:00"
$ws1.Range("K24").Value = "2025-11-03T10:55:36.701838+00:00"
$ws1.Range("K25").Value = "2025-11-03T10:55:36.701847+00:00"
$ws1.Range("K26").Value = "2025-11-03T10:55:36.701855+00:00"
$ws1.Range("K27").Value = "2025-11-03T10:55:43.780492+00:00"
$ws1.Range("K28").Value = "2025-11-03T10:55:43.780524+00:00"
$ws1.Range("K29").Value = "2025-11-03T10:55:43.780560+00:00"
$ws1.Range("K30").Value = "2025-11-03T10:55:43.780589+00:00"
$ws1.Range("K31").Value = "2025-11-03T10:55:43.780609+00:00"
$ws1.Range("K32").Value = "2025-11-03T10:55:45.812428+00:00"
$ws1.Range("K33").Value = "2025-11-03T10:55:45.812445+00:00"
$ws1.Range("K34").Value = "2025-11-03T10:55:45.812453+00:00"
$ws1.Range("K35").Value = "2025-11-03T10:55:48.397722+00:00"
$ws1.Range("K36").Value = "2025-11-03T10:55:48.397752+00:00"
$ws1.Range("K37").Value = "2025-11-03T10:55:48.397772+00:00"
$ws1.Range("K38").Value = "2025-11-03T10:55:48.397793+00:00"
$ws1.Range("K39").Value = "2025-11-03T10:55:48.397811+00:00"
$ws1.Range("K40").Value = "2025-11-03T10:55:48.397826+00:00"
$ws1.Range("K41").Value = "2025-11-03T10:55:48.397866+00:00"
$ws1.Range("K42").Value = "2025-11-03T10:55:48.397900+00:00"
$ws1.Range("K43").Value = "2025-11-03T10:55:50.975196+00:00"
$ws1.Range("K44").Value = "2025-11-03T10:55:50.975230+00:00"
$ws1.Range("K45").Value = "2025-11-03T10:55:56.197421+00:00"
$ws1.Range("K46").Value = "2025-11-03T10:55:58.242309+00:00"
$ws1.Range("K47").Value = "2025-11-03T10:55:58.242343+00:00"
$ws1.Range("K48").Value = "2025-11-03T10:55:58.242365+00:00"

# Remove the two data rows from the "new_injured" sheet, leaving only the header row
$ws3 = $wb.Worksheets.Item("new_injured")
$ws3.Range("A2:G3").EntireRow.Delete()
